# Applies the "Updated cryptos list" refresh: new prices/volumes for most rows,
# plus a handful of rows whose Coin/Link identity shifted (rank reshuffle).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row updates: Row -> @{ Col = NewValue; ... } (only columns that actually changed)
$updates = @{
    2 = @{ "D"="61.066.92"; "E"="  +8.14%  " }
    3 = @{ "D"="3.355.24"; "E"="  +4.20%  " }
    4 = @{ "E"="  -0.07%  " }
    5 = @{ "D"="415.23"; "E"="  +6.16%  " }
    6 = @{ "D"="115.01"; "E"="  +8.07%  " }
    7 = @{ "D"="0.589"; "E"="  +4.79%  " }
    8 = @{ "D"="1.00"; "E"="  -0.01%  " }
    9 = @{ "D"="0.643"; "E"="  +5.24%  " }
    10 = @{ "D"="40.34"; "E"="  +4.51%  " }
    11 = @{ "D"="0.101"; "E"="  +5.88%  " }
    12 = @{ "E"="  +1.33%  " }
    13 = @{ "D"="3.887.21"; "E"="  +4.31%  " }
    14 = @{ "D"="8.47"; "E"="  +5.13%  " }
    15 = @{ "D"="19.93"; "E"="  +6.03%  " }
    16 = @{ "D"="3.383.83"; "E"="  +5.30%  " }
    17 = @{ "E"="  +2.65%  " }
    18 = @{ "D"="60.846.95"; "E"="  +7.94%  " }
    19 = @{ "D"="10.84"; "E"="  +1.11%  " }
    20 = @{ "D"="3.41"; "E"="  +4.10%  " }
    21 = @{ "E"="  +8.19%  " }
    22 = @{ "D"="13.09"; "E"="  +2.64%  " }
    23 = @{ "D"="303.81"; "E"="  +2.63%  " }
    24 = @{ "D"="75.24"; "E"="  +2.90%  " }
    25 = @{ "D"="3.22"; "E"="  +3.90%  " }
    26 = @{ "D"="28.82"; "E"="  +4.12%  " }
    27 = @{ "D"="4.48"; "E"="  +2.42%  " }
    28 = @{ "B"="Filecoin"; "C"="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; "D"="8.01"; "E"="  +4.32%  " }
    29 = @{ "B"="RenderToken"; "C"="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; "D"="7.71"; "E"="  +8.12%  " }
    30 = @{ "B"="Kaspa"; "C"="https://coinranking.com/coin/V8GxkwWow+kaspa-kas"; "D"="0.179"; "E"="  +7.16%  " }
    31 = @{ "B"="Hedera"; "C"="https://coinranking.com/coin/jad286TjB+hedera-hbar"; "D"="0.115"; "E"="  +7.04%  " }
    32 = @{ "B"="Toncoin"; "C"="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; "D"="2.60"; "E"="  +23.48%  " }
    33 = @{ "D"="11.53"; "E"="  +5.93%  " }
    34 = @{ "D"="0.999"; "E"="  -0.02%  " }
    35 = @{ "D"="40.04"; "E"="  +7.94%  " }
    36 = @{ "D"="0.0512"; "E"="  +6.55%  " }
    37 = @{ "D"="52.53"; "E"="  +2.61%  " }
    38 = @{ "D"="3.11"; "E"="  +3.09%  " }
    39 = @{ "E"="  +0.19%  " }
    40 = @{ "D"="3.42"; "E"="  -1.73%  " }
    41 = @{ "D"="137.28"; "E"="  +3.17%  " }
    42 = @{ "B"="Stellar"; "C"="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; "D"="0.124"; "E"="  +4.09%  " }
    43 = @{ "B"="TheGraph"; "C"="https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"; "D"="0.295"; "E"="  +4.55%  " }
    44 = @{ "D"="1.93"; "E"="  +2.54%  " }
    45 = @{ "D"="4.01"; "E"="  +3.32%  " }
    46 = @{ "D"="17.00"; "E"="  +1.16%  " }
    47 = @{ "B"="EnergySwap"; "C"="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; "D"="22.61"; "E"="  +4.36%  " }
    48 = @{ "B"="WEMIXToken"; "C"="https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"; "D"="2.25"; "E"="  +9.41%  " }
    49 = @{ "D"="2.173.08"; "E"="  +2.38%  " }
    50 = @{ "D"="2.40"; "E"="  +2.01%  " }
    51 = @{ "D"="1.98"; "E"="  -1.40%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        $newValue = $rowData[$col]
        if ($col -eq "D") {
            # Price column holds text like "1.00" / "0.999" / "61.066.92" -- force
            # text storage so Excel does not coerce it into a trimmed numeric value,
            # then drop back to the sheet default style (no custom format was applied
            # to these cells originally).
            $cell.NumberFormat = "@"
            $cell.Value = $newValue
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newValue
        }
    }
}
